$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B header (row 1) ---
$ws.Range("B1").Value = "ln_real_wage"

# --- Data rows 2-30: labels (col A) and values (col B) ---
$ws.Range("A2").Value = "group_1"
$ws.Range("B2").Value = 2.6254132011217872
$ws.Range("A3").Value = "group_2"
$ws.Range("B3").Value = 2.0991224233920756
$ws.Range("A4").Value = "difference"
$ws.Range("B4").Value = 0.52629077772971167
$ws.Range("A5").Value = "endowments"
$ws.Range("B5").Value = -0.14134725004044157
$ws.Range("A6").Value = "coefficients"
$ws.Range("B6").Value = 0.36006751649960878
$ws.Range("A7").Value = "interaction"
$ws.Range("B7").Value = 0.30757051127054447
$ws.Range("A8").Value = "LTHS"
$ws.Range("B8").Value = -0.016966524969175206
$ws.Range("A9").Value = "some_college"
$ws.Range("B9").Value = 0.0082961461771727373
$ws.Range("A10").Value = "college"
$ws.Range("B10").Value = 0.0034823329632576927
$ws.Range("A11").Value = "high_school"
$ws.Range("B11").Value = -0.014197203619435208
$ws.Range("A12").Value = "advanced_degree"
$ws.Range("B12").Value = 0.0016387449238859729
$ws.Range("A13").Value = "foreign_born"
$ws.Range("B13").Value = -0.061800372758073749
$ws.Range("A14").Value = "native"
$ws.Range("B14").Value = -0.061800372758073749
$ws.Range("A15").Value = "LTHS"
$ws.Range("B15").Value = -0.31308678434359039
$ws.Range("A16").Value = "some_college"
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = "college"
$ws.Range("B17").Value = 0
$ws.Range("A18").Value = "high_school"
$ws.Range("B18").Value = -0.020999136776596788
$ws.Range("A19").Value = "advanced_degree"
$ws.Range("B19").Value = 0
$ws.Range("A20").Value = "foreign_born"
$ws.Range("B20").Value = 0.019703207601276584
$ws.Range("A21").Value = "native"
$ws.Range("B21").Value = -0.10836764180702121
$ws.Range("A22").Value = "Intercept"
$ws.Range("B22").Value = 0.78281787182554075
$ws.Range("A23").Value = "LTHS"
$ws.Range("B23").Value = 0.2065030291221813
$ws.Range("A24").Value = "some_college"
$ws.Range("B24").Value = -0.036428013282573529
$ws.Range("A25").Value = "college"
$ws.Range("B25").Value = 0.013373256452303629
$ws.Range("A26").Value = "high_school"
$ws.Range("B26").Value = -0.0096581022132280273
$ws.Range("A27").Value = "advanced_degree"
$ws.Range("B27").Value = 0.025307092605966166
$ws.Range("A28").Value = "foreign_born"
$ws.Range("B28").Value = 0.054236624292947413
$ws.Range("A29").Value = "native"
$ws.Range("B29").Value = 0.054236624292947413
$ws.Range("A30").Value = "N"
$ws.Range("B30").Value = 666

# --- Number formats ---
# Column B decimal rows (2-29): 6-decimal custom format
$ws.Range("B2:B29").NumberFormat = "0.######"
# Final row (N count): integer format
$ws.Range("B30").NumberFormat = "0"

# --- Alignment ---
$ws.Range("A2:A30").HorizontalAlignment = -4131
$ws.Range("B2:B30").HorizontalAlignment = -4152
$ws.Range("B1").HorizontalAlignment = -4108

# --- Borders: clear stale border left on what used to be the last row (24), ---
# --- then draw the bottom border under the new last row (30) ---
$ws.Range("A24:B24").Borders.Item(9).LineStyle = -4142
$ws.Range("A30:B30").Borders.Item(9).LineStyle = 1

# --- Update used range / dimension ---
$ws.Range("A1:B30").Select()
